$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (this shifts batsman..sr from D..I to F..K)
$ws.Range("D:E").Insert()

# New header cells for ownTeam / oppTeam
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Populate ownTeam / oppTeam for existing row 2
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Sunrisers Hyderabad"

# Force numeric-looking stat columns to stay text (avoid float re-interpretation)
$statCols = "G2:K2,G3:K3,G4:K4"
$ws.Range($statCols).NumberFormat = "@"

# Add new row 3 data
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 01 2020"
$ws.Range("C3").Value = "Mumbai won by 48 runs"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Mumbai Indians"
$ws.Range("F3").Value = "Ravi Bishnoi "
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "5"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "20.00"

# Add new row 4 data
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " October 08 2020"
$ws.Range("C4").Value = "Sunrisers won by 69 runs"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Ravi Bishnoi "
$ws.Range("G4").Value = "6"
$ws.Range("H4").Value = "7"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "85.71"

# Re-assert the pre-existing row 2 stat values stay text "0" (already text, but keep consistent)
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "-"
